$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "562 WILLIAM BERCZY BLVD MARKHAM ON L6C2P7"
$ws.Range("A2").Value = "MARKHAM WILLIAM BERCZY BLVD ON L6C2P7 562"
$ws.Range("A3").Value = "WILLIAM BERCZY BLVD 562 ON L6C2P7 MARKHAM"
$ws.Range("A4").Value = "WILLIAM BERCZY BLVD MARKHAM ON L6C2P7 562"
$ws.Range("A5").Value = "562 MARKHAM ON L6C2P7 WILLIAM BERCZY BLVD"
$ws.Range("A6").Value = "MARKHAM 562 ON L6C2P7 WILLIAM BERCZY BLVD"

$ws.Range("J5").Select()
